# Apply data updates to the 'Inscricoes' sheet (Table1), columns E (Inscritos),
# F (Pagos), H (Inscricoes homologadas) for specific rows, per commit 'Data update using git'.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 107
$ws.Range("F2").Value = 75
$ws.Range("H2").Value = 81

$ws.Range("E5").Value = 151
$ws.Range("F5").Value = 104
$ws.Range("H5").Value = 115

$ws.Range("E6").Value = 46
$ws.Range("F6").Value = 34
$ws.Range("H6").Value = 44

$ws.Range("E7").Value = 37

$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 5
$ws.Range("H8").Value = 8

$ws.Range("E10").Value = 634
$ws.Range("F10").Value = 342
$ws.Range("H10").Value = 438

$ws.Range("E11").Value = 419

$ws.Range("E12").Value = 635
$ws.Range("F12").Value = 376
$ws.Range("H12").Value = 462

$ws.Range("E13").Value = 153
$ws.Range("F13").Value = 84
$ws.Range("H13").Value = 118

$ws.Range("E16").Value = 221

$ws.Range("E17").Value = 114

$ws.Range("E18").Value = 54

$ws.Range("E19").Value = 16

$ws.Range("E20").Value = 97
$ws.Range("F20").Value = 37
$ws.Range("H20").Value = 74

$ws.Range("E22").Value = 182

$ws.Range("E23").Value = 214
$ws.Range("F23").Value = 108
$ws.Range("H23").Value = 159

$ws.Range("E24").Value = 240
$ws.Range("F24").Value = 137
$ws.Range("H24").Value = 167

$ws.Range("E25").Value = 303
$ws.Range("F25").Value = 165
$ws.Range("H25").Value = 225

$ws.Range("E26").Value = 175
$ws.Range("F26").Value = 106
$ws.Range("H26").Value = 131

$ws.Range("E27").Value = 360
$ws.Range("F27").Value = 192
$ws.Range("H27").Value = 273

$ws.Range("E28").Value = 215

$ws.Range("E29").Value = 183

$ws.Range("E30").Value = 238
$ws.Range("F30").Value = 146
$ws.Range("H30").Value = 198

$ws.Range("E31").Value = 79
$ws.Range("F31").Value = 34
$ws.Range("H31").Value = 62

$ws.Range("E32").Value = 198
$ws.Range("F32").Value = 126
$ws.Range("H32").Value = 164

$ws.Range("E33").Value = 316
$ws.Range("F33").Value = 172
$ws.Range("H33").Value = 262

$ws.Range("E34").Value = 239

$ws.Range("E35").Value = 166
$ws.Range("F35").Value = 113
$ws.Range("H35").Value = 140

$ws.Range("E36").Value = 86
$ws.Range("F36").Value = 54
$ws.Range("H36").Value = 64

$ws.Range("E39").Value = 189

$ws.Range("E40").Value = 286
$ws.Range("F40").Value = 141
$ws.Range("H40").Value = 221

$ws.Range("E41").Value = 421

$ws.Range("E42").Value = 423

$ws.Range("E43").Value = 134

$ws.Range("E44").Value = 341
$ws.Range("F44").Value = 177
$ws.Range("H44").Value = 245

$ws.Range("E45").Value = 163

$ws.Range("E46").Value = 365

$ws.Range("E47").Value = 508
$ws.Range("F47").Value = 276
$ws.Range("H47").Value = 368

$ws.Range("E48").Value = 247
$ws.Range("F48").Value = 112
$ws.Range("H48").Value = 156

$ws.Range("E49").Value = 319
$ws.Range("F49").Value = 156
$ws.Range("H49").Value = 243

$ws.Range("E50").Value = 264
$ws.Range("F50").Value = 139
$ws.Range("H50").Value = 210

$ws.Range("E51").Value = 257
$ws.Range("F51").Value = 125
$ws.Range("H51").Value = 199
